# Refresh cryptocurrency price/volume snapshot per the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.723.06"
$ws.Cells.Item(2, 5).Value = "  +0.37%  "
$ws.Cells.Item(3, 4).Value = "1.600.32"
$ws.Cells.Item(3, 5).Value = "  +0.23%  "
$ws.Cells.Item(4, 5).Value = "  +0.33%  "
$savedStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.54"
$ws.Cells.Item(5, 4).Style = $savedStyle
$ws.Cells.Item(5, 5).Value = "  -0.07%  "
$savedStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.511"
$ws.Cells.Item(6, 4).Style = $savedStyle
$ws.Cells.Item(6, 5).Value = "  -0.86%  "
$ws.Cells.Item(7, 5).Value = "  +0.32%  "
$ws.Cells.Item(8, 5).Value = "  +0.23%  "
$savedStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.248"
$ws.Cells.Item(9, 4).Style = $savedStyle
$ws.Cells.Item(9, 5).Value = "  +0.61%  "
$savedStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.55"
$ws.Cells.Item(10, 4).Style = $savedStyle
$ws.Cells.Item(10, 5).Value = "  +0.08%  "
$ws.Cells.Item(11, 5).Value = "  +0.47%  "
$ws.Cells.Item(12, 4).Value = "1.824.68"
$ws.Cells.Item(12, 5).Value = "  +0.19%  "
$ws.Cells.Item(13, 4).Value = "1.628.02"
$ws.Cells.Item(13, 5).Value = "  +1.82%  "
$ws.Cells.Item(14, 5).Value = "  +0.69%  "
$savedStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.523"
$ws.Cells.Item(15, 4).Style = $savedStyle
$ws.Cells.Item(15, 5).Value = "  +0.26%  "
$savedStyle = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "65.36"
$ws.Cells.Item(16, 4).Style = $savedStyle
$ws.Cells.Item(16, 5).Value = "  +1.53%  "
$ws.Cells.Item(17, 4).Value = "26.690.83"
$ws.Cells.Item(18, 4).Value = "0.0₃0753"
$ws.Cells.Item(18, 5).Value = "  +2.89%  "
$ws.Cells.Item(19, 5).Value = "  +3.98%  "
$ws.Cells.Item(20, 5).Value = "  +0.39%  "
$savedStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "209.17"
$ws.Cells.Item(21, 4).Style = $savedStyle
$ws.Cells.Item(21, 5).Value = "  +0.27%  "
$ws.Cells.Item(22, 5).Value = "  +0.61%  "
$ws.Cells.Item(23, 5).Value = "  +0.77%  "
$savedStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.94"
$ws.Cells.Item(24, 4).Style = $savedStyle
$ws.Cells.Item(24, 5).Value = "  +0.61%  "
$savedStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "142.31"
$ws.Cells.Item(25, 4).Style = $savedStyle
$ws.Cells.Item(25, 5).Value = "  -1.87%  "
$ws.Cells.Item(26, 5).Value = "  +0.20%  "
$savedStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.11"
$ws.Cells.Item(27, 4).Style = $savedStyle
$ws.Cells.Item(27, 5).Value = "  -0.51%  "
$savedStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.115"
$ws.Cells.Item(28, 4).Style = $savedStyle
$ws.Cells.Item(28, 5).Value = "  -0.08%  "
$ws.Cells.Item(29, 5).Value = "  +0.73%  "
$savedStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0523"
$ws.Cells.Item(30, 4).Style = $savedStyle
$ws.Cells.Item(30, 5).Value = "  +3.09%  "
$ws.Cells.Item(31, 5).Value = "  -0.34%  "
$ws.Cells.Item(32, 5).Value = "  +0.84%  "
$savedStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.98"
$ws.Cells.Item(33, 4).Style = $savedStyle
$ws.Cells.Item(33, 5).Value = "  +1.82%  "
$ws.Cells.Item(34, 4).Value = "1.293.45"
$ws.Cells.Item(34, 5).Value = "  +1.27%  "
$ws.Cells.Item(35, 5).Value = "  -5.30%  "
$ws.Cells.Item(36, 5).Value = "  +0.95%  "
$ws.Cells.Item(37, 5).Value = "  +0.41%  "
$ws.Cells.Item(38, 5).Value = "  -0.10%  "
$ws.Cells.Item(39, 5).Value = "  +19.84%  "
$savedStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.826"
$ws.Cells.Item(40, 4).Style = $savedStyle
$ws.Cells.Item(40, 5).Value = "  -2.13%  "
$savedStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.43"
$ws.Cells.Item(41, 4).Style = $savedStyle
$ws.Cells.Item(41, 5).Value = "  -0.93%  "
$ws.Cells.Item(42, 5).Value = "  -0.30%  "
$savedStyle = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.784"
$ws.Cells.Item(43, 4).Style = $savedStyle
$ws.Cells.Item(43, 5).Value = "  -0.20%  "
$savedStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "63.24"
$ws.Cells.Item(44, 4).Style = $savedStyle
$ws.Cells.Item(44, 5).Value = "  -1.86%  "
$ws.Cells.Item(45, 4).Value = "1.736.39"
$ws.Cells.Item(45, 5).Value = "  +0.14%  "
$savedStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "91.33"
$ws.Cells.Item(46, 4).Style = $savedStyle
$ws.Cells.Item(46, 5).Value = "  +1.57%  "
$savedStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.58"
$ws.Cells.Item(47, 4).Style = $savedStyle
$ws.Cells.Item(47, 5).Value = "  -1.69%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0104"
$ws.Cells.Item(48, 5).Value = "  -1.86%  "
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$savedStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.101"
$ws.Cells.Item(49, 4).Style = $savedStyle
$ws.Cells.Item(49, 5).Value = "  -1.29%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$savedStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0510"
$ws.Cells.Item(50, 4).Style = $savedStyle
$ws.Cells.Item(50, 5).Value = "  +0.60%  "
$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$savedStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Cells.Item(51, 4).Style = $savedStyle
$ws.Cells.Item(51, 5).Value = "  +0.33%  "
